$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line Data")

# M8, M9 and M10 now branch to both M11 and M12 instead of just M12
$ws.Range("E9").Value = '["M11", "M12"]'
$ws.Range("E10").Value = '["M11", "M12"]'
$ws.Range("E11").Value = '["M11", "M12"]'

# M11 now links directly to END instead of M12
$ws.Range("E12").Value = "END"

# Row 13 keeps the M12 machine name (now its own distinct entry)
$ws.Range("A13").Value = "M12"

# Header row was resized slightly
$ws.Rows.Item(1).RowHeight = 18.75
